$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 0.02080121282546088
$ws.Range("A2").Value = -0.0099999990691621576
$ws.Range("A3").Value = -0.0089999990810056829
$ws.Range("A4").Value = 0.28398477809229306
$ws.Range("A5").Value = -0.0059999991098553807
$ws.Range("A6").Value = -0.0059999990844907813
$ws.Range("A7").Value = -0.019999998938290631
$ws.Range("A8").Value = -0.019999998934460805
$ws.Range("A9").Value = -0.0059999990769776801
$ws.Range("A10").Value = -0.0059999990761738786
$ws.Range("A11").Value = -0.0044999990916174681
$ws.Range("A12").Value = 0.081974786783316045
$ws.Range("A13").Value = -0.0059999990530705816
$ws.Range("A14").Value = -0.011999998983499793
$ws.Range("A15").Value = -0.0059999990391794711
$ws.Range("A16").Value = -0.034603001594259108
$ws.Range("A17").Value = -0.0059999990284076432
$ws.Range("A18").Value = -0.0089999989967237681
$ws.Range("A19").Value = -0.0089999990788438566
$ws.Range("A20").Value = -0.01115191184000075
$ws.Range("A21").Value = -0.0089999990693927501
$ws.Range("A22").Value = -0.0089999990686502329
$ws.Range("A23").Value = -0.0089999990741089775
$ws.Range("A24").Value = -0.041999998718291742
$ws.Range("A25").Value = -0.041999998711357733
$ws.Range("A26").Value = -0.0059999990813786042
$ws.Range("A27").Value = -0.0059999990772832135
$ws.Range("A28").Value = -0.0059999990589609808
$ws.Range("A29").Value = -0.011999998984823179
$ws.Range("A30").Value = -0.019999998896340632
$ws.Range("A31").Value = -0.0029878909661693598
$ws.Range("A32").Value = -0.029664752671738626
$ws.Range("A33").Value = -0.00599999903044246
